# This script replicates a "regenerate merged AHB files" edit:
# For a set of rows in the single worksheet, the "ÄNDERUNG" marker
# (shared string, stored in column L) is cleared out and the row's
# formatting is switched from the "unmarked / plain" style set
# (s=5 for most cells, s=7 for the L "ÄNDERUNG" cell) to the
# "template / header" style set already used on row 2 of the sheet
# (s=2 for most cells, s=3 for the bold B/M cells, s=4 for the now
# empty L cell). Rows that are the first row of a new numbered
# group (column B value differs from the row above) get the style
# applied across the whole row A:V; all other affected rows only
# get the L cell's style/content changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where the entire row A:V is restyled (first row of a new group)
$fullRows = @(60,63,67,71,75,89,94,98,105,119,140,143,147,152,155,159,161,165,169,173,177,183,186,189,192,195)

# Rows where only column L is restyled (value cleared, style 7 -> 4)
$lOnlyRows = @(61,62,64,65,66,68,69,70,72,73,74,76,77,78,79,80,81,82,83,84,85,86,87,88,90,91,92,93,95,96,97,99,100,101,102,103,104,106,107,108,109,110,111,112,113,114,115,116,117,118,120,121,122,123,141,144,145,146,148,149,150,151,153,154,156,157,158,160,162,163,164,166,167,168,170,171,172,174,175,176,178,179,180,181,182,184,185,187,188,190,191,193,194,196,197)

$xlPasteFormats = -4122

# Template ranges already carrying the desired target formatting (row 2)
$templateRow = $ws.Range("A2:V2")
$templateL = $ws.Range("L2")

foreach ($r in $fullRows) {
    $templateRow.Copy() | Out-Null
    $dst = $ws.Range("A" + $r + ":V" + $r)
    $dst.PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range("L" + $r).ClearContents() | Out-Null
}

foreach ($r in $lOnlyRows) {
    $templateL.Copy() | Out-Null
    $dst = $ws.Range("L" + $r)
    $dst.PasteSpecial($xlPasteFormats) | Out-Null
    $dst.ClearContents() | Out-Null
}

$excel.CutCopyMode = 0
